$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.291.30"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.189.08"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "255.19"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.20%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.628"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "68.26"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.62%  "
$ws.Range("E8").Value = "  -0.04%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.571"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.05%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "58.87"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "36.97"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -5.68%  "
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  +3.21%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "2.512.14"
$ws.Range("E15").Value = "  -1.61%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.868"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.79%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.38"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "2.201.04"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "41.176.36"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  -0.19%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.15"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -0.49%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "233.01"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -3.77%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "11.77"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +19.26%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.86"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +6.16%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.50"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.56%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "169.16"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "20.63"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0744"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.45"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "26.11"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +7.70%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.17"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0297"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.29%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "12.15"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +12.24%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.65"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.91"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "60.84"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -9.13%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.198"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.17%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.62"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.94%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.101"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.21"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +9.63%  "
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.17"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.21"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -9.31%  "
